$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header columns AD (Wins), AE (Losses), AF (Ties) ---
# Copy the existing header style (bold, bordered, centered) from AC1 onto the
# three new header cells before writing their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Fill in the season record for every player row (2-45) ---
$lastRow = 45
$ws.Range("AD2:AD$lastRow").Value = 68
$ws.Range("AE2:AE$lastRow").Value = 94
$ws.Range("AF2:AF$lastRow").Value = 0

Write-Output "season record columns added"
